$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Bonus: Login Page" entry (currently row 20) needs to move down to row 21,
# leaving a new blank row 20 as a separator (mirrors the gap pattern already
# used between rows 13/15 and 17/19). Inserting a row above it shifts it down
# without disturbing anything else.
$ws.Rows.Item(20).Insert()

# Row 16: rename the task.
$ws.Cells.Item(16, 1).Value = "Create Pages & Gut Old App Stuff"

# Row 17: rename the task and flesh out the rest of the row (Hours, amount,
# and a reference link) - set the link (column D) first so the shared-string
# table ends up in the same order as the target workbook.
$ws.Cells.Item(17, 4).Value = "https://www.cssmatic.com/noise-texture"
$ws.Cells.Item(17, 1).Value = "Change Bootstrap & Custom CSS"
$ws.Cells.Item(17, 2).Value = "Hours"
$ws.Cells.Item(17, 3).Value = 1

# Row 19 (new): another task entry, following the same blank-row-gap pattern.
$ws.Cells.Item(19, 1).Value = "Finish Pages And Display Data Nicely"
$ws.Cells.Item(19, 2).Value = "Hours"

# Update the active selection to reflect where editing left off.
$ws.Range("D19").Select()
